# Insert a new weekly record (row 221) into the Brócoli price series.
# All existing rows from 221 downward shift to 222..305; the former last
# row (304) becomes row 305. The workbook's dimension grows from
# A1:R304 to A1:R305.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 221..304 down by one row, creating a blank row 221.
$ws.Rows.Item(221).Insert()

# Populate the new row 221 with the inserted record.
$ws.Range("A221").Value = 7
$ws.Range("B221").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C221").Value = "Ñuble"
$ws.Range("D221").Value = 44795
$ws.Range("E221").Value = 16
$ws.Range("F221").Value = 100112023
$ws.Range("G221").Value = "Brócoli"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Segunda"
$ws.Range("J221").Value = 80
$ws.Range("K221").Value = 900
$ws.Range("L221").Value = 900
$ws.Range("M221").Value = 900
$ws.Range("N221").Value = "$/unidad"
$ws.Range("O221").Value = "Región del Maule"
$ws.Range("P221").Value = 900
$ws.Range("Q221").Value = 1
$ws.Range("R221").Value = "Hortaliza"
